$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1634.375
$ws.Range("I11").Value = 1634.375
$ws.Range("K11").Value = 1634.375
$ws.Range("M11").Value = -1494.375
$ws.Range("H17").Value = 878.1429000000001
$ws.Range("J17").Value = 878.1429000000001
$ws.Range("L17").Value = 2634.4287
$ws.Range("N17").Value = -2970.4287
$ws.Range("H33").Value = 426.88235
$ws.Range("I33").Value = 473.8
$ws.Range("K33").Value = 473.8
$ws.Range("M33").Value = -244.8
$ws.Range("H40").Value = 5712.577
$ws.Range("I40").Value = 3416.25
$ws.Range("J40").Value = 7680.857
$ws.Range("K40").Value = 3416.25
$ws.Range("L40").Value = 7680.857
$ws.Range("M40").Value = -3241.25
$ws.Range("N40").Value = -8030.857
$ws.Range("H41").Value = 16668468
$ws.Range("I41").Value = 722.625
$ws.Range("K41").Value = 722.625
$ws.Range("M41").Value = -282.625
$ws.Range("H51").Value = 2976.2666
$ws.Range("I51").Value = 2198.75
$ws.Range("J51").Value = 3259
$ws.Range("K51").Value = 2198.75
$ws.Range("L51").Value = 3259
$ws.Range("M51").Value = -1714.75
$ws.Range("N51").Value = -4227
$ws.Range("H62").Value = 7816874.5
$ws.Range("I62").Value = 13891580
$ws.Range("K62").Value = 13891580
$ws.Range("M62").Value = -13890956
$ws.Range("H65").Value = 7816874.5
$ws.Range("I65").Value = 13891580
$ws.Range("K65").Value = 69457900
$ws.Range("M65").Value = -69454780
$ws.Range("H74").Value = 6717.1875
$ws.Range("I74").Value = 5605.357
$ws.Range("K74").Value = 5605.357
$ws.Range("M74").Value = -4669.357
$ws.Range("H76").Value = 41713070
$ws.Range("I76").Value = 94771.82000000001
$ws.Range("J76").Value = 76928560
$ws.Range("K76").Value = 94771.82000000001
$ws.Range("L76").Value = 76928560
$ws.Range("M76").Value = -94456.82000000001
$ws.Range("N76").Value = -76929190
$ws.Range("H77").Value = 6717.1875
$ws.Range("I77").Value = 5605.357
$ws.Range("K77").Value = 28026.785
$ws.Range("M77").Value = -23346.785
$ws.Range("H79").Value = 41713070
$ws.Range("I79").Value = 94771.82000000001
$ws.Range("J79").Value = 76928560
$ws.Range("K79").Value = 94771.82000000001
$ws.Range("L79").Value = 76928560
$ws.Range("M79").Value = -93679.82000000001
$ws.Range("N79").Value = -76930744
$ws.Range("H96").Value = 866.75
$ws.Range("I96").Value = 491.33334
$ws.Range("J96").Value = 1993
$ws.Range("K96").Value = 1474.00002
$ws.Range("L96").Value = 5979
$ws.Range("M96").Value = -101.0000199999999
$ws.Range("N96").Value = -8725
$ws.Range("H98").Value = 2550.0667
$ws.Range("I98").Value = 1947.2632
$ws.Range("J98").Value = 3591.2727
$ws.Range("K98").Value = 1947.2632
$ws.Range("L98").Value = 3591.2727
$ws.Range("M98").Value = -449.2632000000001
$ws.Range("N98").Value = -6587.2727
$ws.Range("H106").Value = 1566.9333
$ws.Range("I106").Value = 1188
$ws.Range("K106").Value = 1188
$ws.Range("M106").Value = -557
$ws.Range("H116").Value = 7558.364
$ws.Range("I116").Value = 8254.888999999999
$ws.Range("K116").Value = 8254.888999999999
$ws.Range("M116").Value = -4812.888999999999
$ws.Range("H122").Value = 2550.0667
$ws.Range("I122").Value = 1947.2632
$ws.Range("J122").Value = 3591.2727
$ws.Range("K122").Value = 5841.7896
$ws.Range("L122").Value = 10773.8181
$ws.Range("M122").Value = -3391.7896
$ws.Range("N122").Value = -15673.8181
$ws.Range("H125").Value = 5147.643
$ws.Range("I125").Value = 4983.625
$ws.Range("K125").Value = 44852.625
$ws.Range("M125").Value = -42392.625
$ws.Range("H132").Value = 4782.4473
$ws.Range("I132").Value = 5094.121
$ws.Range("J132").Value = 2725.4
$ws.Range("K132").Value = 15282.363
$ws.Range("L132").Value = 8176.200000000001
$ws.Range("M132").Value = -12752.363
$ws.Range("N132").Value = -13236.2
$ws.Range("H133").Value = 49998.5
$ws.Range("J133").Value = 49998.5
$ws.Range("L133").Value = 49998.5
$ws.Range("N133").Value = -60118.5
$ws.Range("H135").Value = 1608.7046
$ws.Range("I135").Value = 1570.3429
$ws.Range("J135").Value = 1757.8889
$ws.Range("K135").Value = 14133.0861
$ws.Range("L135").Value = 15821.0001
$ws.Range("M135").Value = -11598.0861
$ws.Range("N135").Value = -20891.0001
$ws.Range("H137").Value = 14658.406
$ws.Range("I137").Value = 1493.7693
$ws.Range("K137").Value = 4481.3079
$ws.Range("M137").Value = -1931.3079
$ws.Range("H138").Value = 5100.098
$ws.Range("I138").Value = 2339.8333
$ws.Range("J138").Value = 6605.697
$ws.Range("K138").Value = 7019.499899999999
$ws.Range("L138").Value = 19817.091
$ws.Range("M138").Value = -1879.499899999999
$ws.Range("N138").Value = -30097.091
$ws.Range("H141").Value = 3222.4167
$ws.Range("I141").Value = 3048.3333
$ws.Range("J141").Value = 3744.6667
$ws.Range("K141").Value = 9144.999899999999
$ws.Range("L141").Value = 11234.0001
$ws.Range("M141").Value = -3964.999899999999
$ws.Range("N141").Value = -21594.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 46439.41
$ws.Range("I2").Value = 72387
$ws.Range("J2").Value = 1031.125
$ws.Range("K2").Value = 72387
$ws.Range("L2").Value = 1031.125
$ws.Range("M2").Value = -72274
$ws.Range("N2").Value = -1257.125
$ws.Range("H32").Value = 1368.9354
$ws.Range("I32").Value = 1421.2456
$ws.Range("K32").Value = 1421.2456
$ws.Range("M32").Value = -1134.2456
$ws.Range("H45").Value = 40160.96
$ws.Range("I45").Value = 64011.688
$ws.Range("K45").Value = 64011.688
$ws.Range("M45").Value = -63634.688
$ws.Range("H61").Value = 3862.524
$ws.Range("I61").Value = 3727.2632
$ws.Range("J61").Value = 5147.5
$ws.Range("K61").Value = 3727.2632
$ws.Range("L61").Value = 5147.5
$ws.Range("M61").Value = -3515.2632
$ws.Range("N61").Value = -5571.5
$ws.Range("H63").Value = 8560.5
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = $null
$ws.Range("H66").Value = 8560.5
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = $null
$ws.Range("H69").Value = 188000
$ws.Range("J69").Value = 188000
$ws.Range("L69").Value = 188000
$ws.Range("N69").Value = -189498
$ws.Range("H72").Value = 188000
$ws.Range("J72").Value = 188000
$ws.Range("L72").Value = 564000
$ws.Range("N72").Value = -571488
$ws.Range("H74").Value = 341102.97
$ws.Range("I74").Value = 371583.75
$ws.Range("J74").Value = 66776
$ws.Range("K74").Value = 371583.75
$ws.Range("L74").Value = 66776
$ws.Range("M74").Value = -370709.75
$ws.Range("N74").Value = -68524
$ws.Range("H77").Value = 341102.97
$ws.Range("I77").Value = 371583.75
$ws.Range("J77").Value = 66776
$ws.Range("K77").Value = 1857918.75
$ws.Range("L77").Value = 333880
$ws.Range("M77").Value = -1853550.75
$ws.Range("N77").Value = -342616
$ws.Range("H88").Value = 2961.0625
$ws.Range("I88").Value = 4100.143
$ws.Range("J88").Value = 2075.111
$ws.Range("K88").Value = 4100.143
$ws.Range("L88").Value = 2075.111
$ws.Range("M88").Value = -3694.143
$ws.Range("N88").Value = -2887.111
$ws.Range("H91").Value = 2961.0625
$ws.Range("I91").Value = 4100.143
$ws.Range("J91").Value = 2075.111
$ws.Range("K91").Value = 4100.143
$ws.Range("L91").Value = 2075.111
$ws.Range("M91").Value = -2696.143
$ws.Range("N91").Value = -4883.111
$ws.Range("H110").Value = 1541
$ws.Range("I110").Value = 1469.2
$ws.Range("K110").Value = 1469.2
$ws.Range("M110").Value = 575.8
$ws.Range("H116").Value = 46439.41
$ws.Range("I116").Value = 72387
$ws.Range("J116").Value = 1031.125
$ws.Range("K116").Value = 72387
$ws.Range("L116").Value = 1031.125
$ws.Range("M116").Value = -70093
$ws.Range("N116").Value = -5619.125
$ws.Range("H122").Value = 5771.9
$ws.Range("I122").Value = 4829.1763
$ws.Range("J122").Value = 7004.6924
$ws.Range("K122").Value = 14487.5289
$ws.Range("L122").Value = 21014.0772
$ws.Range("M122").Value = -12037.5289
$ws.Range("N122").Value = -25914.0772
$ws.Range("H132").Value = 226048.2
$ws.Range("I132").Value = 350713.38
$ws.Range("K132").Value = 1052140.14
$ws.Range("M132").Value = -1049610.14
$ws.Range("H136").Value = 3862.524
$ws.Range("I136").Value = 3727.2632
$ws.Range("J136").Value = 5147.5
$ws.Range("K136").Value = 11181.7896
$ws.Range("L136").Value = 15442.5
$ws.Range("M136").Value = -8631.7896
$ws.Range("N136").Value = -20542.5
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null
$ws.Range("H141").Value = 58000
$ws.Range("J141").Value = 58000
$ws.Range("L141").Value = 58000
$ws.Range("N141").Value = -68360

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 45750
$ws.Range("J2").Value = 52666.668
$ws.Range("L2").Value = 52666.668
$ws.Range("N2").Value = -52892.668
$ws.Range("H3").Value = 46439.41
$ws.Range("I3").Value = 72387
$ws.Range("J3").Value = 1031.125
$ws.Range("K3").Value = 72387
$ws.Range("L3").Value = 1031.125
$ws.Range("M3").Value = -72273
$ws.Range("N3").Value = -1259.125
$ws.Range("H8").Value = 2875
$ws.Range("I8").Value = 2875
$ws.Range("K8").Value = 2875
$ws.Range("M8").Value = -2735
$ws.Range("H86").Value = 1701800.1
$ws.Range("I86").Value = 2834616.8
$ws.Range("J86").Value = 2575
$ws.Range("K86").Value = 2834616.8
$ws.Range("L86").Value = 2575
$ws.Range("M86").Value = -2833493.8
$ws.Range("N86").Value = -4821
$ws.Range("H89").Value = 1701800.1
$ws.Range("I89").Value = 2834616.8
$ws.Range("J89").Value = 2575
$ws.Range("K89").Value = 14173084
$ws.Range("L89").Value = 12875
$ws.Range("M89").Value = -14167468
$ws.Range("N89").Value = -24107
$ws.Range("H94").Value = 1463
$ws.Range("I94").Value = 1253
$ws.Range("K94").Value = 1253
$ws.Range("M94").Value = -802
$ws.Range("H105").Value = 2450.7
$ws.Range("I105").Value = 2450.7
$ws.Range("K105").Value = 2450.7
$ws.Range("M105").Value = -703.6999999999998
$ws.Range("H134").Value = 28067.094
$ws.Range("I134").Value = 1399.2667
$ws.Range("J134").Value = 89608.234
$ws.Range("K134").Value = 4197.800099999999
$ws.Range("L134").Value = 268824.702
$ws.Range("M134").Value = -1662.800099999999
$ws.Range("N134").Value = -273894.702

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1962.375
$ws.Range("I16").Value = 1528.4286
$ws.Range("K16").Value = 1528.4286
$ws.Range("M16").Value = -1241.4286
$ws.Range("H31").Value = 621562.5600000001
$ws.Range("I31").Value = 1542257.5
$ws.Range("K31").Value = 1542257.5
$ws.Range("M31").Value = -1541962.5
$ws.Range("H34").Value = 621562.5600000001
$ws.Range("I34").Value = 1542257.5
$ws.Range("K34").Value = 1542257.5
$ws.Range("M34").Value = -1542055.5
$ws.Range("H58").Value = 10081.348
$ws.Range("I58").Value = 4705.5557
$ws.Range("K58").Value = 4705.5557
$ws.Range("M58").Value = -4502.5557
$ws.Range("H94").Value = 818.75
$ws.Range("J94").Value = 1004.625
$ws.Range("L94").Value = 1004.625
$ws.Range("N94").Value = -1906.625
$ws.Range("H99").Value = 6854.364
$ws.Range("I99").Value = 5099.6
$ws.Range("J99").Value = 8316.666999999999
$ws.Range("K99").Value = 5099.6
$ws.Range("L99").Value = 8316.666999999999
$ws.Range("M99").Value = -3601.6
$ws.Range("N99").Value = -11312.667
$ws.Range("H100").Value = 46657.332
$ws.Range("J100").Value = 46657.332
$ws.Range("L100").Value = 46657.332
$ws.Range("N100").Value = -48821.332
$ws.Range("H107").Value = 1237.5652
$ws.Range("I107").Value = 817
$ws.Range("K107").Value = 817
$ws.Range("M107").Value = 1103
$ws.Range("H113").Value = 1962.375
$ws.Range("I113").Value = 1528.4286
$ws.Range("K113").Value = 1528.4286
$ws.Range("M113").Value = 641.5714
$ws.Range("H122").Value = 2831.8823
$ws.Range("I122").Value = 2250.4
$ws.Range("J122").Value = 3662.5715
$ws.Range("K122").Value = 6751.200000000001
$ws.Range("L122").Value = 10987.7145
$ws.Range("M122").Value = -4301.200000000001
$ws.Range("N122").Value = -15887.7145
$ws.Range("H126").Value = 6854.364
$ws.Range("I126").Value = 5099.6
$ws.Range("J126").Value = 8316.666999999999
$ws.Range("K126").Value = 15298.8
$ws.Range("L126").Value = 24950.001
$ws.Range("M126").Value = -12828.8
$ws.Range("N126").Value = -29890.001
$ws.Range("H132").Value = 3567.2222
$ws.Range("I132").Value = 2571.3
$ws.Range("J132").Value = 4812.125
$ws.Range("K132").Value = 7713.900000000001
$ws.Range("L132").Value = 14436.375
$ws.Range("M132").Value = -5183.900000000001
$ws.Range("N132").Value = -19496.375
$ws.Range("H134").Value = 462888
$ws.Range("I134").Value = 3251.7856
$ws.Range("K134").Value = 9755.356800000001
$ws.Range("M134").Value = -7220.356800000001
$ws.Range("H136").Value = 10081.348
$ws.Range("I136").Value = 4705.5557
$ws.Range("K136").Value = 14116.6671
$ws.Range("M136").Value = -11566.6671
$ws.Range("H141").Value = 64973.918
$ws.Range("J141").Value = 64973.918
$ws.Range("L141").Value = 64973.918
$ws.Range("N141").Value = -75333.91800000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 4110
$ws.Range("J2").Value = 1833.3334
$ws.Range("L2").Value = 11000.0004
$ws.Range("N2").Value = -11226.0004
$ws.Range("H18").Value = 607.25
$ws.Range("I18").Value = 354
$ws.Range("J18").Value = 1029.3334
$ws.Range("K18").Value = 1062
$ws.Range("L18").Value = 3088.0002
$ws.Range("M18").Value = -893
$ws.Range("N18").Value = -3426.0002
$ws.Range("H37").Value = 116461.875
$ws.Range("J37").Value = 116461.875
$ws.Range("L37").Value = 349385.625
$ws.Range("N37").Value = -349609.625
$ws.Range("H40").Value = 155
$ws.Range("I40").Value = 103.333336
$ws.Range("J40").Value = 387.5
$ws.Range("K40").Value = 413.333344
$ws.Range("L40").Value = 1550
$ws.Range("M40").Value = -344.333344
$ws.Range("N40").Value = -1688
$ws.Range("H47").Value = 8876
$ws.Range("I47").Value = 5425.6
$ws.Range("K47").Value = 16276.8
$ws.Range("M47").Value = -15845.8
$ws.Range("H62").Value = 4299.8335
$ws.Range("J62").Value = 4333
$ws.Range("L62").Value = 12999
$ws.Range("N62").Value = -14371
$ws.Range("H65").Value = 4299.8335
$ws.Range("J65").Value = 4333
$ws.Range("L65").Value = 38997
$ws.Range("N65").Value = -45861
$ws.Range("H68").Value = 1874.7142
$ws.Range("I68").Value = 1875
$ws.Range("J68").Value = 1874.6923
$ws.Range("K68").Value = 5625
$ws.Range("L68").Value = 5624.0769
$ws.Range("M68").Value = -4814
$ws.Range("N68").Value = -7246.0769
$ws.Range("H69").Value = 2999.875
$ws.Range("I69").Value = 2999.5
$ws.Range("J69").Value = 3000
$ws.Range("K69").Value = 8998.5
$ws.Range("L69").Value = 9000
$ws.Range("M69").Value = -8187.5
$ws.Range("N69").Value = -10622
$ws.Range("H71").Value = 1874.7142
$ws.Range("I71").Value = 1875
$ws.Range("J71").Value = 1874.6923
$ws.Range("K71").Value = 16875
$ws.Range("L71").Value = 16872.2307
$ws.Range("M71").Value = -12819
$ws.Range("N71").Value = -24984.2307
$ws.Range("H72").Value = 2999.875
$ws.Range("I72").Value = 2999.5
$ws.Range("J72").Value = 3000
$ws.Range("K72").Value = 26995.5
$ws.Range("L72").Value = 27000
$ws.Range("M72").Value = -22939.5
$ws.Range("N72").Value = -35112
$ws.Range("H76").Value = 4507.5
$ws.Range("I76").Value = 5000
$ws.Range("K76").Value = 15000
$ws.Range("M76").Value = -14617
$ws.Range("H79").Value = 4507.5
$ws.Range("I79").Value = 5000
$ws.Range("K79").Value = 15000
$ws.Range("M79").Value = -13674
$ws.Range("H80").Value = 2015
$ws.Range("J80").Value = 2250
$ws.Range("L80").Value = 6750
$ws.Range("N80").Value = -8622
$ws.Range("H83").Value = 2015
$ws.Range("J83").Value = 2250
$ws.Range("L83").Value = 20250
$ws.Range("N83").Value = -29610
$ws.Range("H107").Value = 24105.652
$ws.Range("J107").Value = 28304.46
$ws.Range("L107").Value = 84913.38
$ws.Range("N107").Value = -88753.38
$ws.Range("H112").Value = 170114.33
$ws.Range("I112").Value = 202977.2
$ws.Range("J112").Value = 5800
$ws.Range("K112").Value = 608931.6000000001
$ws.Range("L112").Value = 17400
$ws.Range("M112").Value = -607823.6000000001
$ws.Range("N112").Value = -19616
$ws.Range("H113").Value = 6748742.5
$ws.Range("I113").Value = 7856865.5
$ws.Range("J113").Value = 100004
$ws.Range("K113").Value = 23570596.5
$ws.Range("L113").Value = 300012
$ws.Range("M113").Value = -23568426.5
$ws.Range("N113").Value = -304352
$ws.Range("H117").Value = 900
$ws.Range("J117").Value = 900
$ws.Range("L117").Value = 2700
$ws.Range("N117").Value = -9584
$ws.Range("H122").Value = 611443.1
$ws.Range("I122").Value = 2198242
$ws.Range("J122").Value = 1135.8462
$ws.Range("K122").Value = 19784178
$ws.Range("L122").Value = 10222.6158
$ws.Range("M122").Value = -19781728
$ws.Range("N122").Value = -15122.6158
$ws.Range("H128").Value = 148399.6
$ws.Range("I128").Value = 148399.6
$ws.Range("K128").Value = 445198.8
$ws.Range("M128").Value = -440218.8
$ws.Range("H129").Value = 1148
$ws.Range("J129").Value = 1997
$ws.Range("L129").Value = 5991
$ws.Range("N129").Value = -15991
$ws.Range("H130").Value = 249.5
$ws.Range("J130").Value = 249
$ws.Range("L130").Value = 747
$ws.Range("N130").Value = -10787
$ws.Range("H132").Value = 2543624
$ws.Range("I132").Value = 9092587
$ws.Range("J132").Value = 24792.385
$ws.Range("K132").Value = 81833283
$ws.Range("L132").Value = 223131.465
$ws.Range("M132").Value = -81830753
$ws.Range("N132").Value = -228191.465
$ws.Range("H141").Value = 3687.375
$ws.Range("I141").Value = 2800
$ws.Range("K141").Value = 8400
$ws.Range("M141").Value = -3220

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9460.799999999999
$ws.Range("I70").Value = 9460.799999999999
$ws.Range("K70").Value = 9460.799999999999
$ws.Range("M70").Value = -9190.799999999999
$ws.Range("H73").Value = 9460.799999999999
$ws.Range("I73").Value = 9460.799999999999
$ws.Range("K73").Value = 9460.799999999999
$ws.Range("M73").Value = -8524.799999999999
$ws.Range("H80").Value = 913452.5
$ws.Range("I80").Value = 669443.5600000001
$ws.Range("J80").Value = 1436328.9
$ws.Range("K80").Value = 669443.5600000001
$ws.Range("L80").Value = 1436328.9
$ws.Range("M80").Value = -668445.5600000001
$ws.Range("N80").Value = -1438324.9
$ws.Range("H83").Value = 913452.5
$ws.Range("I83").Value = 669443.5600000001
$ws.Range("J83").Value = 1436328.9
$ws.Range("K83").Value = 3347217.8
$ws.Range("L83").Value = 7181644.5
$ws.Range("M83").Value = -3342225.8
$ws.Range("N83").Value = -7191628.5
$ws.Range("H102").Value = 41979.89
$ws.Range("I102").Value = 115944
$ws.Range("J102").Value = 4997.8335
$ws.Range("K102").Value = 115944
$ws.Range("L102").Value = 4997.8335
$ws.Range("M102").Value = -114322
$ws.Range("N102").Value = -8241.833500000001
$ws.Range("H118").Value = 50000
$ws.Range("J118").Value = 50000
$ws.Range("L118").Value = 50000
$ws.Range("N118").Value = -53314
$ws.Range("H122").Value = 617684.6
$ws.Range("I122").Value = 739719
$ws.Range("J122").Value = 7512.6665
$ws.Range("K122").Value = 2219157
$ws.Range("L122").Value = 22537.9995
$ws.Range("M122").Value = -2216707
$ws.Range("N122").Value = -27437.9995
$ws.Range("H132").Value = 67809.766
$ws.Range("I132").Value = 22462.938
$ws.Range("K132").Value = 67388.814
$ws.Range("M132").Value = -64858.814
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 353613.28
$ws.Range("I7").Value = 593089.75
$ws.Range("J7").Value = 14354.917
$ws.Range("K7").Value = 593089.75
$ws.Range("L7").Value = 14354.917
$ws.Range("M7").Value = -592977.75
$ws.Range("N7").Value = -14578.917
$ws.Range("H16").Value = 2461.25
$ws.Range("I16").Value = 2461.25
$ws.Range("K16").Value = 2461.25
$ws.Range("M16").Value = -2291.25
$ws.Range("H22").Value = 4000.3333
$ws.Range("I22").Value = 4000.3333
$ws.Range("K22").Value = 4000.3333
$ws.Range("M22").Value = -3705.3333
$ws.Range("H27").Value = 4000.3333
$ws.Range("I27").Value = 4000.3333
$ws.Range("K27").Value = 4000.3333
$ws.Range("M27").Value = -3893.3333
$ws.Range("H40").Value = 422684.5
$ws.Range("I40").Value = 506066.75
$ws.Range("J40").Value = 5773.25
$ws.Range("K40").Value = 506066.75
$ws.Range("L40").Value = 5773.25
$ws.Range("M40").Value = -505930.75
$ws.Range("N40").Value = -6045.25
$ws.Range("H46").Value = 2810.8462
$ws.Range("I46").Value = 2373.238
$ws.Range("J46").Value = 3321.389
$ws.Range("K46").Value = 2373.238
$ws.Range("L46").Value = 3321.389
$ws.Range("M46").Value = -2185.238
$ws.Range("N46").Value = -3697.389
$ws.Range("H50").Value = 15000000
$ws.Range("I50").Value = 15000000
$ws.Range("K50").Value = 15000000
$ws.Range("M50").Value = -14999363
$ws.Range("H53").Value = 63331.668
$ws.Range("I53").Value = 79997.5
$ws.Range("J53").Value = 30000
$ws.Range("K53").Value = 79997.5
$ws.Range("L53").Value = 30000
$ws.Range("M53").Value = -79479.5
$ws.Range("N53").Value = -31036
$ws.Range("H55").Value = 34483256
$ws.Range("I55").Value = 228.4375
$ws.Range("J55").Value = 76923900
$ws.Range("K55").Value = 228.4375
$ws.Range("L55").Value = 76923900
$ws.Range("M55").Value = -55.4375
$ws.Range("N55").Value = -76924246
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null
$ws.Range("H68").Value = 94681.17999999999
$ws.Range("I68").Value = 3966.6667
$ws.Range("J68").Value = 128699.125
$ws.Range("K68").Value = 3966.6667
$ws.Range("L68").Value = 128699.125
$ws.Range("M68").Value = -3217.6667
$ws.Range("N68").Value = -130197.125
$ws.Range("H71").Value = 94681.17999999999
$ws.Range("I71").Value = 3966.6667
$ws.Range("J71").Value = 128699.125
$ws.Range("K71").Value = 19833.3335
$ws.Range("L71").Value = 643495.625
$ws.Range("M71").Value = -16089.3335
$ws.Range("N71").Value = -650983.625
$ws.Range("H74").Value = 34924.25
$ws.Range("H77").Value = 34924.25
$ws.Range("H93").Value = 2544.4285
$ws.Range("I93").Value = 3052.4
$ws.Range("K93").Value = 3052.4
$ws.Range("M93").Value = -1804.4
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null
$ws.Range("H122").Value = 922849
$ws.Range("I122").Value = 5395.6665
$ws.Range("K122").Value = 16186.9995
$ws.Range("M122").Value = -13736.9995
$ws.Range("H126").Value = 353613.28
$ws.Range("I126").Value = 593089.75
$ws.Range("J126").Value = 14354.917
$ws.Range("K126").Value = 1779269.25
$ws.Range("L126").Value = 43064.751
$ws.Range("M126").Value = -1776799.25
$ws.Range("N126").Value = -48004.751
$ws.Range("H132").Value = 6090.107
$ws.Range("I132").Value = 4948.6313
$ws.Range("J132").Value = 8499.888999999999
$ws.Range("K132").Value = 14845.8939
$ws.Range("L132").Value = 25499.667
$ws.Range("M132").Value = -12315.8939
$ws.Range("N132").Value = -30559.667
$ws.Range("H136").Value = 1695706.8
$ws.Range("I136").Value = 2859599.5
$ws.Range("K136").Value = 8578798.5
$ws.Range("M136").Value = -8576248.5
$ws.Range("H138").Value = 75000
$ws.Range("J138").Value = 75000
$ws.Range("L138").Value = 75000
$ws.Range("N138").Value = -85280

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 500010000
$ws.Range("I26").Value = 20000
$ws.Range("K26").Value = 20000
$ws.Range("M26").Value = -19707
$ws.Range("H55").Value = 9082.166999999999
$ws.Range("J55").Value = 9746.5
$ws.Range("L55").Value = 9746.5
$ws.Range("N55").Value = -10300.5
$ws.Range("H62").Value = 7291
$ws.Range("I62").Value = 6582.1665
$ws.Range("J62").Value = 7999.8335
$ws.Range("K62").Value = 6582.1665
$ws.Range("L62").Value = 7999.8335
$ws.Range("M62").Value = -5958.1665
$ws.Range("N62").Value = -9247.833500000001
$ws.Range("H65").Value = 7291
$ws.Range("I65").Value = 6582.1665
$ws.Range("J65").Value = 7999.8335
$ws.Range("K65").Value = 32910.8325
$ws.Range("L65").Value = 39999.1675
$ws.Range("M65").Value = -29790.8325
$ws.Range("N65").Value = -46239.1675
$ws.Range("H75").Value = 36500
$ws.Range("J75").Value = 36500
$ws.Range("L75").Value = 36500
$ws.Range("N75").Value = -38372
$ws.Range("H78").Value = 36500
$ws.Range("J78").Value = 36500
$ws.Range("L78").Value = 109500
$ws.Range("N78").Value = -118860
$ws.Range("H81").Value = 2959.5386
$ws.Range("I81").Value = 2153.1428
$ws.Range("K81").Value = 4306.2856
$ws.Range("M81").Value = -3245.2856
$ws.Range("H84").Value = 2959.5386
$ws.Range("I84").Value = 2153.1428
$ws.Range("K84").Value = 21531.428
$ws.Range("M84").Value = -16227.428
$ws.Range("H100").Value = 2187
$ws.Range("I100").Value = 2356.5715
$ws.Range("K100").Value = 4713.143
$ws.Range("M100").Value = -4172.143
$ws.Range("H107").Value = 67681.47
$ws.Range("I107").Value = 100871.9
$ws.Range("J107").Value = 1300.6
$ws.Range("K107").Value = 302615.7
$ws.Range("L107").Value = 3901.8
$ws.Range("M107").Value = -300695.7
$ws.Range("N107").Value = -7741.799999999999
$ws.Range("H113").Value = 1687.6666
$ws.Range("I113").Value = 1667.8667
$ws.Range("K113").Value = 5003.6001
$ws.Range("M113").Value = -2833.6001
$ws.Range("H122").Value = 3480
$ws.Range("I122").Value = 2908.375
$ws.Range("K122").Value = 8725.125
$ws.Range("M122").Value = -6275.125
$ws.Range("H132").Value = 21724.725
$ws.Range("I132").Value = 1412.975
$ws.Range("K132").Value = 4238.924999999999
$ws.Range("M132").Value = -1708.924999999999
$ws.Range("H136").Value = 421202.66
$ws.Range("I136").Value = 420419.75
$ws.Range("K136").Value = 1261259.25
$ws.Range("M136").Value = -1258709.25
$ws.Range("H140").Value = 120868.625
$ws.Range("J140").Value = 120868.625
$ws.Range("L140").Value = 120868.625
$ws.Range("N140").Value = -131228.625
$ws.Range("H141").Value = 52142.855
$ws.Range("J141").Value = 52142.855
$ws.Range("L141").Value = 52142.855
$ws.Range("N141").Value = -62502.855

Write-Host "Done applying changes"